# Update the AddEmployee data table with new sample values and move the
# active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddEmployee")

# Header row (row 1) stays the same: FirstName | MiddleName | LastName

# Replace the sample data rows (A2:C5) with the new values.
$ws.Range("A2").Value = "John"
$ws.Range("B2").Value = "K"
$ws.Range("C2").Value = "Doe"

$ws.Range("A3").Value = "Katie"
$ws.Range("B3").Value = "K"
$ws.Range("C3").Value = "Ball"

$ws.Range("A4").Value = "Donald"
$ws.Range("B4").Value = "K"
$ws.Range("C4").Value = "Trump"

$ws.Range("A5").Value = "Mohammed"
$ws.Range("B5").Value = "K"
$ws.Range("C5").Value = "Salah"

# Move the active selection to E11, matching the saved view state.
$ws.Range("E11").Select()
